$d = $word.ActiveDocument

$replacements = @(
    @("939÷6=156, 3", "918÷2=459, 0"),
    @("219÷6=36, 3", "990÷3=330, 0"),
    @("871÷5=174, 1", "922÷3=307, 1"),
    @("895÷8=111, 7", "835÷2=417, 1"),
    @("881÷5=176, 1", "170÷7=24, 2"),
    @("235÷6=39, 1", "980÷7=140, 0"),
    @("954÷8=119, 2", "561÷5=112, 1"),
    @("546÷2=273, 0", "402÷6=67, 0"),
    @("266÷6=44, 2", "391÷4=97, 3"),
    @("400÷7=57, 1", "501÷9=55, 6"),
    @("573÷6=95, 3", "584÷6=97, 2"),
    @("940÷5=188, 0", "319÷5=63, 4"),
    @("975÷5=195, 0", "450÷2=225, 0"),
    @("277÷6=46, 1", "823÷9=91, 4"),
    @("292÷9=32, 4", "701÷2=350, 1"),
    @("514÷5=102, 4", "134÷5=26, 4"),
    @("867÷5=173, 2", "631÷8=78, 7"),
    @("325÷6=54, 1", "111÷8=13, 7"),
    @("641÷4=160, 1", "962÷4=240, 2"),
    @("219÷7=31, 2", "151÷7=21, 4"),
    @("232÷4=58, 0", "710÷4=177, 2"),
    @("593÷7=84, 5", "492÷2=246, 0"),
    @("815÷5=163, 0", "550÷7=78, 4"),
    @("921÷6=153, 3", "326÷7=46, 4"),
    @("217÷2=108, 1", "790÷9=87, 7")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

$d.Save()
